# Applies targeted cell updates to multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# matching the upstream scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 77059.766
$ws.Range("I11").Value = 77059.766
$ws.Range("K11").Value = 77059.766
$ws.Range("M11").Value = -76919.766

$ws.Range("H63").Value = 10000
$ws.Range("J63").Value = 10000
$ws.Range("L63").Value = 10000
$ws.Range("N63").Value = -11248

$ws.Range("H66").Value = 10000
$ws.Range("J66").Value = 10000
$ws.Range("L66").Value = 30000
$ws.Range("N66").Value = -36240

$ws.Range("H86").Value = 4202.905
$ws.Range("I86").Value = 1647.5
$ws.Range("J86").Value = 5225.067
$ws.Range("K86").Value = 1647.5
$ws.Range("L86").Value = 5225.067
$ws.Range("M86").Value = -524.5
$ws.Range("N86").Value = -7471.067

$ws.Range("H89").Value = 4202.905
$ws.Range("I89").Value = 1647.5
$ws.Range("J89").Value = 5225.067
$ws.Range("K89").Value = 8237.5
$ws.Range("L89").Value = 26125.335
$ws.Range("M89").Value = -2621.5
$ws.Range("N89").Value = -37357.335

$ws.Range("H136").Value = 59663.332
$ws.Range("J136").Value = 59663.332
$ws.Range("L136").Value = 59663.332
$ws.Range("N136").Value = -69863.33199999999

$ws.Range("H137").Value = 1276.9048
$ws.Range("I137").Value = 1108.875
$ws.Range("K137").Value = 3326.625
$ws.Range("M137").Value = -776.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2876.55
$ws.Range("I32").Value = 1906.5914
$ws.Range("J32").Value = 15763.143
$ws.Range("K32").Value = 1906.5914
$ws.Range("L32").Value = 15763.143
$ws.Range("M32").Value = -1619.5914
$ws.Range("N32").Value = -16337.143

$ws.Range("H61").Value = 1570.5846
$ws.Range("I61").Value = 1030.762
$ws.Range("J61").Value = 2556.348
$ws.Range("K61").Value = 1030.762
$ws.Range("L61").Value = 2556.348
$ws.Range("M61").Value = -818.7619999999999
$ws.Range("N61").Value = -2980.348

$ws.Range("H132").Value = 1531.4857
$ws.Range("I132").Value = 1226.1482
$ws.Range("J132").Value = 2562
$ws.Range("K132").Value = 3678.4446
$ws.Range("L132").Value = 7686
$ws.Range("M132").Value = -1148.4446
$ws.Range("N132").Value = -12746

$ws.Range("H136").Value = 1570.5846
$ws.Range("I136").Value = 1030.762
$ws.Range("J136").Value = 2556.348
$ws.Range("K136").Value = 3092.286
$ws.Range("L136").Value = 7669.044
$ws.Range("M136").Value = -542.2860000000001
$ws.Range("N136").Value = -12769.044

$ws.Range("H138").Value = 63246.332
$ws.Range("I138").Value = 58900
$ws.Range("J138").Value = 65419.5
$ws.Range("K138").Value = 58900
$ws.Range("L138").Value = 65419.5
$ws.Range("M138").Value = -53760
$ws.Range("N138").Value = -75699.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 64128.562
$ws.Range("I20").Value = 85004.086
$ws.Range("J20").Value = 1502
$ws.Range("K20").Value = 85004.086
$ws.Range("L20").Value = 1502
$ws.Range("M20").Value = -84757.086
$ws.Range("N20").Value = -1996

$ws.Range("H97").Value = 11711.5
$ws.Range("I97").Value = 5615.3335
$ws.Range("J97").Value = 30000
$ws.Range("K97").Value = 5615.3335
$ws.Range("L97").Value = 30000
$ws.Range("M97").Value = -4624.3335
$ws.Range("N97").Value = -31982

$ws.Range("H106").Value = 14134.2
$ws.Range("J106").Value = 14134.2
$ws.Range("L106").Value = 14134.2
$ws.Range("N106").Value = -16658.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1115.5625
$ws.Range("I94").Value = 853
$ws.Range("K94").Value = 853
$ws.Range("M94").Value = -402

$ws.Range("H99").Value = 8959.412
$ws.Range("I99").Value = 2857
$ws.Range("K99").Value = 2857
$ws.Range("M99").Value = -1359

$ws.Range("H126").Value = 8959.412
$ws.Range("I126").Value = 2857
$ws.Range("K126").Value = 8571
$ws.Range("M126").Value = -6101

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 8770.097
$ws.Range("J55").Value = 5905.9653
$ws.Range("L55").Value = 17717.8959
$ws.Range("N55").Value = -18071.8959

$ws.Range("H57").Value = 4472.5
$ws.Range("I57").Value = 4390
$ws.Range("K57").Value = 13170
$ws.Range("M57").Value = -12611

$ws.Range("H88").Value = 7175
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 7175
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 21525
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = -22381

$ws.Range("H91").Value = 7175
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 7175
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 21525
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = -24489

$ws.Range("H130").Value = 1098.3334
$ws.Range("I130").Value = 820
$ws.Range("J130").Value = 1933.3334
$ws.Range("K130").Value = 2460
$ws.Range("L130").Value = 5800.0002
$ws.Range("M130").Value = 2560
$ws.Range("N130").Value = -15840.0002

$ws.Range("H131").Value = 813.55
$ws.Range("J131").Value = 867.5454999999999
$ws.Range("L131").Value = 2602.6365
$ws.Range("N131").Value = -12682.6365

$ws.Range("H136").Value = 4500
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = $null
$ws.Range("N136").Value = -23700

$ws.Range("H138").Value = 3307.5
$ws.Range("I138").Value = 4030
$ws.Range("J138").Value = 3066.6667
$ws.Range("K138").Value = 12090
$ws.Range("L138").Value = 9200.000100000001
$ws.Range("M138").Value = -6950
$ws.Range("N138").Value = -19480.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 23266.666
$ws.Range("J63").Value = 23266.666
$ws.Range("L63").Value = 23266.666
$ws.Range("N63").Value = -24638.666

$ws.Range("H66").Value = 23266.666
$ws.Range("J66").Value = 23266.666
$ws.Range("L66").Value = 69799.99800000001
$ws.Range("N66").Value = -76663.99800000001

$ws.Range("H70").Value = 89669.836
$ws.Range("I70").Value = 130823.125
$ws.Range("J70").Value = 7363.25
$ws.Range("K70").Value = 130823.125
$ws.Range("L70").Value = 7363.25
$ws.Range("M70").Value = -130553.125
$ws.Range("N70").Value = -7903.25

$ws.Range("H73").Value = 89669.836
$ws.Range("I73").Value = 130823.125
$ws.Range("J73").Value = 7363.25
$ws.Range("K73").Value = 130823.125
$ws.Range("L73").Value = 7363.25
$ws.Range("M73").Value = -129887.125
$ws.Range("N73").Value = -9235.25

$ws.Range("H132").Value = 2937
$ws.Range("I132").Value = 2961.4546
$ws.Range("J132").Value = 2903.375
$ws.Range("K132").Value = 8884.363799999999
$ws.Range("L132").Value = 8710.125
$ws.Range("M132").Value = -6354.363799999999
$ws.Range("N132").Value = -13770.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 359.4375
$ws.Range("I55").Value = 203.88889
$ws.Range("J55").Value = 559.4286
$ws.Range("K55").Value = 203.88889
$ws.Range("L55").Value = 559.4286
$ws.Range("M55").Value = -30.88889
$ws.Range("N55").Value = -905.4286

$ws.Range("H132").Value = 2708.4082
$ws.Range("I132").Value = 2808.5676
$ws.Range("J132").Value = 2399.5833
$ws.Range("K132").Value = 8425.702799999999
$ws.Range("L132").Value = 7198.749899999999
$ws.Range("M132").Value = -5895.702799999999
$ws.Range("N132").Value = -12258.7499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 9705.4
$ws.Range("I32").Value = 4818
$ws.Range("J32").Value = 11800
$ws.Range("K32").Value = 4818
$ws.Range("L32").Value = 11800
$ws.Range("M32").Value = -4501
$ws.Range("N32").Value = -12434

$ws.Range("H132").Value = 3722.24
$ws.Range("I132").Value = 6538.125
$ws.Range("J132").Value = 2397.1177
$ws.Range("K132").Value = 19614.375
$ws.Range("L132").Value = 7191.353099999999
$ws.Range("M132").Value = -17084.375
$ws.Range("N132").Value = -12251.3531
